# Pneumatic System Bill of Materials - add eBoot breadboard/jumper-wire kit
# to the "recommended items" list of the BOM. Mirrors the author's manual
# Excel edit: a new row is inserted above the first subtotal (row 6), the
# new row is filled in with the part's details, and the subtotal formulas
# are updated to include it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 6 - everything currently on/after row 6
# (the first subtotal, the "Recommended Items" section, hyperlinks, etc.)
# shifts down by one row automatically, just like pressing Ctrl+Shift+"+"
# in Excel with row 6 selected.
$ws.Rows.Item(6).Insert()

# Row 6 is tall because the note in column H wraps across several lines.
$ws.Rows.Item(6).RowHeight = 72

# Fill in the new row's data. Set the H (Notes) / A (Item) / C (Manufacturer)
# text cells first so the new shared-string entries land in the same order
# as the source edit (Notes, then Item, then Manufacturer).
$ws.Range("H6").Value2 = "At least one breadboard and a couple jumper wires are needed for the circuitry. This pack has more than needed but is cheaper than other packs that have only the amount of resources needed  https://www.amazon.com/eBoot-400-Point-Solderless-Breadboard-Flexible/dp/B071D7V9HD/ref=sr_1_5?ie=UTF8&qid=1544299357&sr=8-5&keywords=breadboard+and+wires "
$ws.Range("A6").Value2 = "eBoot 3 Pieces 400-Point Solderless Circuit Breadboard with 65 Pieces M/M Flexible Breadboard Jumper Wires"
$ws.Range("C6").Value2 = "eBoot"
$ws.Range("B6").Value2 = "Amazon.com"

# Match existing row styling.
$ws.Range("B6").Style = $ws.Range("B5").Style
$ws.Range("C6").Style = $ws.Range("C5").Style
$ws.Range("E6").Style = $ws.Range("E5").Style
$ws.Range("F6").Style = $ws.Range("F5").Style
$ws.Range("G6").Style = $ws.Range("G5").Style
$ws.Range("H6").Style = $ws.Range("H5").Style

# The Item cell (A6) uses the same style as the similar part-number style
# text already used elsewhere in the sheet (dark grey font, wrap text).
$ws.Range("A6").Style = $ws.Range("D5").Style

# Part number / UPC for the kit. It's formatted as a zero-padded 5-digit-plus
# code (numeric "00000" style), like a barcode/UPC lookup field.
$ws.Range("D6").Style = $ws.Range("D5").Style
$ws.Range("D6").NumberFormat = "00000"
$ws.Range("D6").Value2 = 712971918559

$ws.Range("E6").Value2 = 8.99
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Formula = "=F6*E6"

# Update the first subtotal so it includes the newly inserted row.
$ws.Range("G7").Formula = "=SUM(G2:G6)"

$ws.Calculate()

# Match the saved selection from the source edit.
$ws.Range("D4").Select()

$wb.Save()
